# Update countries & provincias Spain
# Applies the "Datos actualizados" timestamp bump, refreshed per-country
# counters, the Malaui/Nicaragua correction, and reshuffles the Chile row
# so it sorts immediately after Pakistan (ahead of Israel/Austria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed "Datos actualizados" timestamp (16:22 -> 16:52)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 16:52"

# Row 4: Estados Unidos (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 1067382
$ws.Range("C4").Value = 3188
$ws.Range("D4").Value = 147480
$ws.Range("E4").Value = 858053
$ws.Range("F4").Value = 18851
$ws.Range("G4").Value = 193
$ws.Range("H4").Value = 61849

# Row 9: Alemania
$ws.Range("B9").Value = 161985
$ws.Range("C9").Value = 446
$ws.Range("D9").Value = 123500
$ws.Range("E9").Value = 31981
$ws.Range("F9").Value = 2415
$ws.Range("G9").Value = 37
$ws.Range("H9").Value = 6504

# Rows 29-31: Chile moves up to sort right after Pakistan, ahead of
# Israel and Austria, each shifting down one row with refreshed Chile data.
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 16023
$ws.Range("C29").Value = 1138
$ws.Range("D29").Value = 8580
$ws.Range("E29").Value = 7216
$ws.Range("F29").Value = 392
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 227

$ws.Range("A30").Value = "Israel"
$ws.Range("B30").Value = 15870
$ws.Range("C30").Value = 36
$ws.Range("D30").Value = 8412
$ws.Range("E30").Value = 7239
$ws.Range("F30").Value = 117
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 219

$ws.Range("A31").Value = "Austria"
$ws.Range("B31").Value = 15452
$ws.Range("C31").Value = 50
$ws.Range("D31").Value = 12907
$ws.Range("E31").Value = 1961
$ws.Range("F31").Value = 128
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 584

# Row 66: Oman - Muertes hoy
$ws.Range("F66").Value = 17

# Row 177: Malaui
$ws.Range("B177").Value = 37
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 27

# Row 198: Nicaragua
$ws.Range("B198").Value = 14
$ws.Range("C198").Value = 1
$ws.Range("E198").Value = 4
